$d = $word.ActiveDocument

# The date line reads "Күні: 20.10.2023ж". The day portion "20" must become
# "23" (Күні: 23.10.2023ж). We anchor on the literal "Күні: " label first so
# that we only touch the intended "20" run and not the "20" inside "2023".
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("Күні: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$dayRange = $d.Range($anchor.End, $anchor.End + 2)
if ($dayRange.Text -eq "20") {
    $dayRange.Text = "23"
}
